# New weekly price record for "Poroto verde" (Macroferia Regional de Talca).
# Insert a new row at row 33 - this shifts the existing data rows 33..119
# down to 34..120 (and extends the used range to A1:R120), exactly like
# Excel's "Insert Sheet Rows" on row 33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(33).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44544
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112031
$ws.Range("G33").Value = "Poroto verde"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 13000
$ws.Range("N33").Value = "`$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 520
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
